# issue #5: stock data from json to db
# The stock sheet ("股票") gains three new columns:
#   - "category" inserted right after "property_category" (new col I)
#   - "source_file" and "index" appended after "legislator_id" (new cols M, N)
# Existing rows are populated with the matching values: category="normal",
# source_file="tmp9f521", index=the row's legislator id number (115 / 116).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # "股票" (stock) sheet

# Insert a new column I ("category"), shifting old date/legislator_name/
# legislator_id (I,J,K) one place to the right (J,K,L). Formatting of the
# inserted column follows the column to its right automatically.
$ws.Columns("I").Insert()

# Extend the table with two more trailing columns (M, N) by copying the
# formatting of the last existing column (L, "legislator_id") so the new
# cells pick up the same cell styles used throughout the sheet.
$ws.Cells.Item(1, 12).Copy($ws.Cells.Item(1, 13))
$ws.Cells.Item(1, 12).Copy($ws.Cells.Item(1, 14))
$ws.Cells.Item(2, 12).Copy($ws.Cells.Item(2, 13))
$ws.Cells.Item(2, 12).Copy($ws.Cells.Item(2, 14))
$ws.Cells.Item(3, 12).Copy($ws.Cells.Item(3, 13))
$ws.Cells.Item(3, 12).Copy($ws.Cells.Item(3, 14))

# Header row (row 1)
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Row 2 data
$ws.Cells.Item(2, 9).Value = "normal"
$ws.Cells.Item(2, 13).Value = "tmp9f521"
$ws.Cells.Item(2, 14).Value = 115

# Row 3 data
$ws.Cells.Item(3, 9).Value = "normal"
$ws.Cells.Item(3, 13).Value = "tmp9f521"
$ws.Cells.Item(3, 14).Value = 116
